# Fruta / hortaliza, semanal
# Insert 3 new weekly records before the existing row 231 (Vega Monumental
# Concepción - Uva), shifting the existing rows 231:250 down to 234:253,
# and fill the 3 new rows with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 231, pushing the old 231:250 block down to 234:253.
$ws.Rows("231:233").Insert()

# --- New row 231: Autumn Royal -----------------------------------------
$ws.Cells.Item(231, 1).Value = 11
$ws.Cells.Item(231, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(231, 3).Value = "Bíobío"
$ws.Cells.Item(231, 4).Value = 45077
$ws.Cells.Item(231, 5).Value = 8
$ws.Cells.Item(231, 6).Value = "Fruta"
$ws.Cells.Item(231, 7).Value = 100109
$ws.Cells.Item(231, 8).Value = "Uva"
$ws.Cells.Item(231, 9).Value = 100109001
$ws.Cells.Item(231, 10).Value = "Uva"
$ws.Cells.Item(231, 11).Value = "Autumn Royal"
$ws.Cells.Item(231, 12).Value = "Primera"
$ws.Cells.Item(231, 13).Value = 170
$ws.Cells.Item(231, 14).Value = 10000
$ws.Cells.Item(231, 15).Value = 11000
$ws.Cells.Item(231, 16).Value = 10529
$ws.Cells.Item(231, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(231, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(231, 19).Value = 585
$ws.Cells.Item(231, 20).Value = 18

# --- New row 232: Crimpson Seedless -------------------------------------
$ws.Cells.Item(232, 1).Value = 11
$ws.Cells.Item(232, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(232, 3).Value = "Bíobío"
$ws.Cells.Item(232, 4).Value = 45077
$ws.Cells.Item(232, 5).Value = 8
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100109
$ws.Cells.Item(232, 8).Value = "Uva"
$ws.Cells.Item(232, 9).Value = 100109001
$ws.Cells.Item(232, 10).Value = "Uva"
$ws.Cells.Item(232, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(232, 12).Value = "Primera"
$ws.Cells.Item(232, 13).Value = 220
$ws.Cells.Item(232, 14).Value = 10000
$ws.Cells.Item(232, 15).Value = 11000
$ws.Cells.Item(232, 16).Value = 10455
$ws.Cells.Item(232, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(232, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(232, 19).Value = 581
$ws.Cells.Item(232, 20).Value = 18

# --- New row 233: Red Globe ---------------------------------------------
$ws.Cells.Item(233, 1).Value = 11
$ws.Cells.Item(233, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(233, 3).Value = "Bíobío"
$ws.Cells.Item(233, 4).Value = 45077
$ws.Cells.Item(233, 5).Value = 8
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100109
$ws.Cells.Item(233, 8).Value = "Uva"
$ws.Cells.Item(233, 9).Value = 100109001
$ws.Cells.Item(233, 10).Value = "Uva"
$ws.Cells.Item(233, 11).Value = "Red Globe"
$ws.Cells.Item(233, 12).Value = "Primera"
$ws.Cells.Item(233, 13).Value = 130
$ws.Cells.Item(233, 14).Value = 10000
$ws.Cells.Item(233, 15).Value = 11000
$ws.Cells.Item(233, 16).Value = 10615
$ws.Cells.Item(233, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(233, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(233, 19).Value = 590
$ws.Cells.Item(233, 20).Value = 18
